# Remove the summary/support sheets that are no longer part of the workbook
$wb = $excel.ActiveWorkbook
foreach ($name in @("Summary", "Department_Summary", "Configuration")) {
    if ($wb.Worksheets | Where-Object { $_.Name -eq $name }) {
        $wb.Worksheets($name).Delete() | Out-Null
    }
}

# Rework the Exam_Schedule sheet: updated day/session slots for the existing
# exams (re-balanced across 2025-11-20 .. 2025-11-28) plus newly scheduled
# exams extending the sheet through 2025-12-05.
$ws = $wb.Worksheets("Exam_Schedule")

# The "date" (H) and "original_preferred" (L) columns hold plain text that
# looks like ISO dates. Force the cells to text first so Excel does not
# auto-convert them into date serial numbers, then drop the now-unneeded
# "@" number format so the cells end up with no special style, same as the
# rest of the plain-text data cells.
$dateTextRangeH = $ws.Range("H2:H40")
$dateTextRangeL = $ws.Range("L2:L40")
$dateTextRangeH.NumberFormat = "@"
$dateTextRangeL.NumberFormat = "@"

# Row 2: CS264 - Data Structures Lab
$ws.Cells.Item(2, 1).Value = "CS264"
$ws.Cells.Item(2, 2).Value = "Data Structures Lab"
$ws.Cells.Item(2, 3).Value = "Lab"
$ws.Cells.Item(2, 4).Value = "2.0 hours"
$ws.Cells.Item(2, 5).Value = 120
$ws.Cells.Item(2, 6).Value = "CSE"
$ws.Cells.Item(2, 7).Value = 3
$ws.Cells.Item(2, 8).Value = "2025-11-20"
$ws.Cells.Item(2, 9).Value = "Thursday"
$ws.Cells.Item(2, 10).Value = "Morning"
$ws.Cells.Item(2, 11).Value = "09:00 - 12:00"
$ws.Cells.Item(2, 12).Value = "2024-12-10"
$ws.Cells.Item(2, 13).Value = "Scheduled"

# Row 3: EC101 - Electronics Theory
$ws.Cells.Item(3, 1).Value = "EC101"
$ws.Cells.Item(3, 2).Value = "Electronics Theory"
$ws.Cells.Item(3, 3).Value = "Theory"
$ws.Cells.Item(3, 4).Value = "3.0 hours"
$ws.Cells.Item(3, 5).Value = 180
$ws.Cells.Item(3, 6).Value = "ECE"
$ws.Cells.Item(3, 7).Value = 1
$ws.Cells.Item(3, 8).Value = "2025-11-20"
$ws.Cells.Item(3, 9).Value = "Thursday"
$ws.Cells.Item(3, 10).Value = "Morning"
$ws.Cells.Item(3, 11).Value = "09:00 - 12:00"
$ws.Cells.Item(3, 12).Value = "2024-12-16"
$ws.Cells.Item(3, 13).Value = "Scheduled"

# Row 4: CS101 - Programming Fundamentals
$ws.Cells.Item(4, 1).Value = "CS101"
$ws.Cells.Item(4, 2).Value = "Programming Fundamentals"
$ws.Cells.Item(4, 3).Value = "Theory"
$ws.Cells.Item(4, 4).Value = "3.0 hours"
$ws.Cells.Item(4, 5).Value = 180
$ws.Cells.Item(4, 6).Value = "CSE"
$ws.Cells.Item(4, 7).Value = 1
$ws.Cells.Item(4, 8).Value = "2025-11-20"
$ws.Cells.Item(4, 9).Value = "Thursday"
$ws.Cells.Item(4, 10).Value = "Afternoon"
$ws.Cells.Item(4, 11).Value = "14:00 - 17:00"
$ws.Cells.Item(4, 12).Value = "2024-12-15"
$ws.Cells.Item(4, 13).Value = "Scheduled"

# Row 5: DA261 - Statistical Programming
$ws.Cells.Item(5, 1).Value = "DA261"
$ws.Cells.Item(5, 2).Value = "Statistical Programming"
$ws.Cells.Item(5, 3).Value = "Lab"
$ws.Cells.Item(5, 4).Value = "2.0 hours"
$ws.Cells.Item(5, 5).Value = 120
$ws.Cells.Item(5, 6).Value = "DSAI"
$ws.Cells.Item(5, 7).Value = 3
$ws.Cells.Item(5, 8).Value = "2025-11-20"
$ws.Cells.Item(5, 9).Value = "Thursday"
$ws.Cells.Item(5, 10).Value = "Afternoon"
$ws.Cells.Item(5, 11).Value = "14:00 - 17:00"
$ws.Cells.Item(5, 12).Value = "2025-04-21"
$ws.Cells.Item(5, 13).Value = "Scheduled"

# Row 6: CS263 - Data Structures
$ws.Cells.Item(6, 1).Value = "CS263"
$ws.Cells.Item(6, 2).Value = "Data Structures"
$ws.Cells.Item(6, 3).Value = "Theory"
$ws.Cells.Item(6, 4).Value = "3.0 hours"
$ws.Cells.Item(6, 5).Value = 180
$ws.Cells.Item(6, 6).Value = "CSE"
$ws.Cells.Item(6, 7).Value = 3
$ws.Cells.Item(6, 8).Value = "2025-11-21"
$ws.Cells.Item(6, 9).Value = "Friday"
$ws.Cells.Item(6, 10).Value = "Morning"
$ws.Cells.Item(6, 11).Value = "09:00 - 12:00"
$ws.Cells.Item(6, 12).Value = "2024-12-17"
$ws.Cells.Item(6, 13).Value = "Scheduled"

# Row 7: DA262 - Data Handling
$ws.Cells.Item(7, 1).Value = "DA262"
$ws.Cells.Item(7, 2).Value = "Data Handling"
$ws.Cells.Item(7, 3).Value = "Theory"
$ws.Cells.Item(7, 4).Value = "3.0 hours"
$ws.Cells.Item(7, 5).Value = 180
$ws.Cells.Item(7, 6).Value = "DSAI"
$ws.Cells.Item(7, 7).Value = 3
$ws.Cells.Item(7, 8).Value = "2025-11-21"
$ws.Cells.Item(7, 9).Value = "Friday"
$ws.Cells.Item(7, 10).Value = "Morning"
$ws.Cells.Item(7, 11).Value = "09:00 - 12:00"
$ws.Cells.Item(7, 12).Value = "2025-04-22"
$ws.Cells.Item(7, 13).Value = "Scheduled"

# Row 8: MA161 - Statistics
$ws.Cells.Item(8, 1).Value = "MA161"
$ws.Cells.Item(8, 2).Value = "Statistics"
$ws.Cells.Item(8, 3).Value = "Theory"
$ws.Cells.Item(8, 4).Value = "2.0 hours"
$ws.Cells.Item(8, 5).Value = 120
$ws.Cells.Item(8, 6).Value = "CSE"
$ws.Cells.Item(8, 7).Value = 1
$ws.Cells.Item(8, 8).Value = "2025-11-21"
$ws.Cells.Item(8, 9).Value = "Friday"
$ws.Cells.Item(8, 10).Value = "Afternoon"
$ws.Cells.Item(8, 11).Value = "14:00 - 17:00"
$ws.Cells.Item(8, 12).Value = "2024-12-18"
$ws.Cells.Item(8, 13).Value = "Scheduled"

# Row 9: CS304 - Artificial Intelligence
$ws.Cells.Item(9, 1).Value = "CS304"
$ws.Cells.Item(9, 2).Value = "Artificial Intelligence"
$ws.Cells.Item(9, 3).Value = "Theory"
$ws.Cells.Item(9, 4).Value = "3.0 hours"
$ws.Cells.Item(9, 5).Value = 180
$ws.Cells.Item(9, 6).Value = "DSAI"
$ws.Cells.Item(9, 7).Value = 3
$ws.Cells.Item(9, 8).Value = "2025-11-21"
$ws.Cells.Item(9, 9).Value = "Friday"
$ws.Cells.Item(9, 10).Value = "Afternoon"
$ws.Cells.Item(9, 11).Value = "14:00 - 17:00"
$ws.Cells.Item(9, 12).Value = "2025-04-23"
$ws.Cells.Item(9, 13).Value = "Scheduled"

# Row 10: DS161 - Introduction to Programming
$ws.Cells.Item(10, 1).Value = "DS161"
$ws.Cells.Item(10, 2).Value = "Introduction to Programming"
$ws.Cells.Item(10, 3).Value = "Theory"
$ws.Cells.Item(10, 4).Value = "2.0 hours"
$ws.Cells.Item(10, 5).Value = 120
$ws.Cells.Item(10, 6).Value = "CSE"
$ws.Cells.Item(10, 7).Value = 1
$ws.Cells.Item(10, 8).Value = "2025-11-24"
$ws.Cells.Item(10, 9).Value = "Monday"
$ws.Cells.Item(10, 10).Value = "Morning"
$ws.Cells.Item(10, 11).Value = "09:00 - 12:00"
$ws.Cells.Item(10, 12).Value = "2024-12-19"
$ws.Cells.Item(10, 13).Value = "Scheduled"

# Row 11: CS307 - Machine Learning
$ws.Cells.Item(11, 1).Value = "CS307"
$ws.Cells.Item(11, 2).Value = "Machine Learning"
$ws.Cells.Item(11, 3).Value = "Theory"
$ws.Cells.Item(11, 4).Value = "3.0 hours"
$ws.Cells.Item(11, 5).Value = 180
$ws.Cells.Item(11, 6).Value = "DSAI"
$ws.Cells.Item(11, 7).Value = 3
$ws.Cells.Item(11, 8).Value = "2025-11-24"
$ws.Cells.Item(11, 9).Value = "Monday"
$ws.Cells.Item(11, 10).Value = "Morning"
$ws.Cells.Item(11, 11).Value = "09:00 - 12:00"
$ws.Cells.Item(11, 12).Value = "2025-04-24"
$ws.Cells.Item(11, 13).Value = "Scheduled"

# Row 12: MA162 - Probability
$ws.Cells.Item(12, 1).Value = "MA162"
$ws.Cells.Item(12, 2).Value = "Probability"
$ws.Cells.Item(12, 3).Value = "Theory"
$ws.Cells.Item(12, 4).Value = "2.0 hours"
$ws.Cells.Item(12, 5).Value = 120
$ws.Cells.Item(12, 6).Value = "CSE"
$ws.Cells.Item(12, 7).Value = 1
$ws.Cells.Item(12, 8).Value = "2025-11-24"
$ws.Cells.Item(12, 9).Value = "Monday"
$ws.Cells.Item(12, 10).Value = "Afternoon"
$ws.Cells.Item(12, 11).Value = "14:00 - 17:00"
$ws.Cells.Item(12, 12).Value = "2024-12-20"
$ws.Cells.Item(12, 13).Value = "Scheduled"

# Row 13: EC301 - Digital Signal Processing
$ws.Cells.Item(13, 1).Value = "EC301"
$ws.Cells.Item(13, 2).Value = "Digital Signal Processing"
$ws.Cells.Item(13, 3).Value = "Theory"
$ws.Cells.Item(13, 4).Value = "3.0 hours"
$ws.Cells.Item(13, 5).Value = 180
$ws.Cells.Item(13, 6).Value = "ECE"
$ws.Cells.Item(13, 7).Value = 3
$ws.Cells.Item(13, 8).Value = "2025-11-24"
$ws.Cells.Item(13, 9).Value = "Monday"
$ws.Cells.Item(13, 10).Value = "Afternoon"
$ws.Cells.Item(13, 11).Value = "14:00 - 17:00"
$ws.Cells.Item(13, 12).Value = "2025-04-25"
$ws.Cells.Item(13, 13).Value = "Scheduled"

# Row 14: EC161 - Digital Design
$ws.Cells.Item(14, 1).Value = "EC161"
$ws.Cells.Item(14, 2).Value = "Digital Design"
$ws.Cells.Item(14, 3).Value = "Theory"
$ws.Cells.Item(14, 4).Value = "2.0 hours"
$ws.Cells.Item(14, 5).Value = 120
$ws.Cells.Item(14, 6).Value = "CSE"
$ws.Cells.Item(14, 7).Value = 1
$ws.Cells.Item(14, 8).Value = "2025-11-25"
$ws.Cells.Item(14, 9).Value = "Tuesday"
$ws.Cells.Item(14, 10).Value = "Morning"
$ws.Cells.Item(14, 11).Value = "09:00 - 12:00"
$ws.Cells.Item(14, 12).Value = "2024-12-21"
$ws.Cells.Item(14, 13).Value = "Scheduled"

# Row 15: HS201 - Happiness & Wellbeing
$ws.Cells.Item(15, 1).Value = "HS201"
$ws.Cells.Item(15, 2).Value = "Happiness & Wellbeing"
$ws.Cells.Item(15, 3).Value = "Theory"
$ws.Cells.Item(15, 4).Value = "3.0 hours"
$ws.Cells.Item(15, 5).Value = 180
$ws.Cells.Item(15, 6).Value = "ECE"
$ws.Cells.Item(15, 7).Value = 3
$ws.Cells.Item(15, 8).Value = "2025-11-25"
$ws.Cells.Item(15, 9).Value = "Tuesday"
$ws.Cells.Item(15, 10).Value = "Morning"
$ws.Cells.Item(15, 11).Value = "09:00 - 12:00"
$ws.Cells.Item(15, 12).Value = "2025-04-28"
$ws.Cells.Item(15, 13).Value = "Scheduled"

# Row 16: CS161 - Problem Solving
$ws.Cells.Item(16, 1).Value = "CS161"
$ws.Cells.Item(16, 2).Value = "Problem Solving"
$ws.Cells.Item(16, 3).Value = "Theory"
$ws.Cells.Item(16, 4).Value = "3.0 hours"
$ws.Cells.Item(16, 5).Value = 180
$ws.Cells.Item(16, 6).Value = "CSE"
$ws.Cells.Item(16, 7).Value = 1
$ws.Cells.Item(16, 8).Value = "2025-11-25"
$ws.Cells.Item(16, 9).Value = "Tuesday"
$ws.Cells.Item(16, 10).Value = "Afternoon"
$ws.Cells.Item(16, 11).Value = "14:00 - 17:00"
$ws.Cells.Item(16, 12).Value = "2024-12-22"
$ws.Cells.Item(16, 13).Value = "Scheduled"

# Row 17: EC302 - Introduction to VLSI Design
$ws.Cells.Item(17, 1).Value = "EC302"
$ws.Cells.Item(17, 2).Value = "Introduction to VLSI Design"
$ws.Cells.Item(17, 3).Value = "Theory"
$ws.Cells.Item(17, 4).Value = "3.0 hours"
$ws.Cells.Item(17, 5).Value = 180
$ws.Cells.Item(17, 6).Value = "ECE"
$ws.Cells.Item(17, 7).Value = 3
$ws.Cells.Item(17, 8).Value = "2025-11-25"
$ws.Cells.Item(17, 9).Value = "Tuesday"
$ws.Cells.Item(17, 10).Value = "Afternoon"
$ws.Cells.Item(17, 11).Value = "14:00 - 17:00"
$ws.Cells.Item(17, 12).Value = "2025-04-26"
$ws.Cells.Item(17, 13).Value = "Scheduled"

# Row 18: HS161 - English Language
$ws.Cells.Item(18, 1).Value = "HS161"
$ws.Cells.Item(18, 2).Value = "English Language"
$ws.Cells.Item(18, 3).Value = "Theory"
$ws.Cells.Item(18, 4).Value = "3.0 hours"
$ws.Cells.Item(18, 5).Value = 180
$ws.Cells.Item(18, 6).Value = "CSE"
$ws.Cells.Item(18, 7).Value = 1
$ws.Cells.Item(18, 8).Value = "2025-11-26"
$ws.Cells.Item(18, 9).Value = "Wednesday"
$ws.Cells.Item(18, 10).Value = "Morning"
$ws.Cells.Item(18, 11).Value = "09:00 - 12:00"
$ws.Cells.Item(18, 12).Value = "2024-12-23"
$ws.Cells.Item(18, 13).Value = "Scheduled"

# Row 19: DS456 - Machine Learning
$ws.Cells.Item(19, 1).Value = "DS456"
$ws.Cells.Item(19, 2).Value = "Machine Learning"
$ws.Cells.Item(19, 3).Value = "Theory"
$ws.Cells.Item(19, 4).Value = "3.0 hours"
$ws.Cells.Item(19, 5).Value = 180
$ws.Cells.Item(19, 6).Value = "DSAI"
$ws.Cells.Item(19, 7).Value = 5
$ws.Cells.Item(19, 8).Value = "2025-11-26"
$ws.Cells.Item(19, 9).Value = "Wednesday"
$ws.Cells.Item(19, 10).Value = "Morning"
$ws.Cells.Item(19, 11).Value = "09:00 - 12:00"
$ws.Cells.Item(19, 12).Value = "2025-05-16"
$ws.Cells.Item(19, 13).Value = "Scheduled"

# Row 20: PH151 - Introduction to Physics
$ws.Cells.Item(20, 1).Value = "PH151"
$ws.Cells.Item(20, 2).Value = "Introduction to Physics"
$ws.Cells.Item(20, 3).Value = "Theory"
$ws.Cells.Item(20, 4).Value = "2.0 hours"
$ws.Cells.Item(20, 5).Value = 120
$ws.Cells.Item(20, 6).Value = "CSE"
$ws.Cells.Item(20, 7).Value = 1
$ws.Cells.Item(20, 8).Value = "2025-11-26"
$ws.Cells.Item(20, 9).Value = "Wednesday"
$ws.Cells.Item(20, 10).Value = "Afternoon"
$ws.Cells.Item(20, 11).Value = "14:00 - 17:00"
$ws.Cells.Item(20, 12).Value = "2024-12-24"
$ws.Cells.Item(20, 13).Value = "Scheduled"

# Row 21: DS302 - Computer Communication
$ws.Cells.Item(21, 1).Value = "DS302"
$ws.Cells.Item(21, 2).Value = "Computer Communication"
$ws.Cells.Item(21, 3).Value = "Theory"
$ws.Cells.Item(21, 4).Value = "3.0 hours"
$ws.Cells.Item(21, 5).Value = 180
$ws.Cells.Item(21, 6).Value = "DSAI"
$ws.Cells.Item(21, 7).Value = 5
$ws.Cells.Item(21, 8).Value = "2025-11-26"
$ws.Cells.Item(21, 9).Value = "Wednesday"
$ws.Cells.Item(21, 10).Value = "Afternoon"
$ws.Cells.Item(21, 11).Value = "14:00 - 17:00"
$ws.Cells.Item(21, 12).Value = "2025-05-21"
$ws.Cells.Item(21, 13).Value = "Scheduled"

# Row 22: CS151 - Introduction to C Programming
$ws.Cells.Item(22, 1).Value = "CS151"
$ws.Cells.Item(22, 2).Value = "Introduction to C Programming"
$ws.Cells.Item(22, 3).Value = "Lab"
$ws.Cells.Item(22, 4).Value = "2.0 hours"
$ws.Cells.Item(22, 5).Value = 120
$ws.Cells.Item(22, 6).Value = "CSE"
$ws.Cells.Item(22, 7).Value = 1
$ws.Cells.Item(22, 8).Value = "2025-11-27"
$ws.Cells.Item(22, 9).Value = "Thursday"
$ws.Cells.Item(22, 10).Value = "Morning"
$ws.Cells.Item(22, 11).Value = "09:00 - 12:00"
$ws.Cells.Item(22, 12).Value = "2024-12-25"
$ws.Cells.Item(22, 13).Value = "Scheduled"

# Row 23: DS303 - Algorithms and Data Structures
$ws.Cells.Item(23, 1).Value = "DS303"
$ws.Cells.Item(23, 2).Value = "Algorithms and Data Structures"
$ws.Cells.Item(23, 3).Value = "Theory"
$ws.Cells.Item(23, 4).Value = "3.0 hours"
$ws.Cells.Item(23, 5).Value = 180
$ws.Cells.Item(23, 6).Value = "DSAI"
$ws.Cells.Item(23, 7).Value = 5
$ws.Cells.Item(23, 8).Value = "2025-11-27"
$ws.Cells.Item(23, 9).Value = "Thursday"
$ws.Cells.Item(23, 10).Value = "Morning"
$ws.Cells.Item(23, 11).Value = "09:00 - 12:00"
$ws.Cells.Item(23, 12).Value = "2025-05-22"
$ws.Cells.Item(23, 13).Value = "Scheduled"

# Row 24: HS157 - Computational Thinking
$ws.Cells.Item(24, 1).Value = "HS157"
$ws.Cells.Item(24, 2).Value = "Computational Thinking"
$ws.Cells.Item(24, 3).Value = "Theory"
$ws.Cells.Item(24, 4).Value = "1.0 hours"
$ws.Cells.Item(24, 5).Value = 60
$ws.Cells.Item(24, 6).Value = "CSE"
$ws.Cells.Item(24, 7).Value = 1
$ws.Cells.Item(24, 8).Value = "2025-11-27"
$ws.Cells.Item(24, 9).Value = "Thursday"
$ws.Cells.Item(24, 10).Value = "Afternoon"
$ws.Cells.Item(24, 11).Value = "14:00 - 17:00"
$ws.Cells.Item(24, 12).Value = "2024-12-26"
$ws.Cells.Item(24, 13).Value = "Scheduled"

# Row 25: EC303 - Random Processes
$ws.Cells.Item(25, 1).Value = "EC303"
$ws.Cells.Item(25, 2).Value = "Random Processes"
$ws.Cells.Item(25, 3).Value = "Theory"
$ws.Cells.Item(25, 4).Value = "2.0 hours"
$ws.Cells.Item(25, 5).Value = 120
$ws.Cells.Item(25, 6).Value = "ECE"
$ws.Cells.Item(25, 7).Value = 5
$ws.Cells.Item(25, 8).Value = "2025-11-27"
$ws.Cells.Item(25, 9).Value = "Thursday"
$ws.Cells.Item(25, 10).Value = "Afternoon"
$ws.Cells.Item(25, 11).Value = "14:00 - 17:00"
$ws.Cells.Item(25, 12).Value = "2025-05-24"
$ws.Cells.Item(25, 13).Value = "Scheduled"

# Row 26: HS156 - Holistic Personality Development
$ws.Cells.Item(26, 1).Value = "HS156"
$ws.Cells.Item(26, 2).Value = "Holistic Personality Development"
$ws.Cells.Item(26, 3).Value = "Theory"
$ws.Cells.Item(26, 4).Value = "1.0 hours"
$ws.Cells.Item(26, 5).Value = 60
$ws.Cells.Item(26, 6).Value = "CSE"
$ws.Cells.Item(26, 7).Value = 1
$ws.Cells.Item(26, 8).Value = "2025-11-28"
$ws.Cells.Item(26, 9).Value = "Friday"
$ws.Cells.Item(26, 10).Value = "Morning"
$ws.Cells.Item(26, 11).Value = "09:00 - 12:00"
$ws.Cells.Item(26, 12).Value = "2024-12-27"
$ws.Cells.Item(26, 13).Value = "Scheduled"

# Row 27: EC304 - Signals & Systems
$ws.Cells.Item(27, 1).Value = "EC304"
$ws.Cells.Item(27, 2).Value = "Signals & Systems"
$ws.Cells.Item(27, 3).Value = "Theory"
$ws.Cells.Item(27, 4).Value = "3.0 hours"
$ws.Cells.Item(27, 5).Value = 180
$ws.Cells.Item(27, 6).Value = "ECE"
$ws.Cells.Item(27, 7).Value = 5
$ws.Cells.Item(27, 8).Value = "2025-11-28"
$ws.Cells.Item(27, 9).Value = "Friday"
$ws.Cells.Item(27, 10).Value = "Morning"
$ws.Cells.Item(27, 11).Value = "09:00 - 12:00"
$ws.Cells.Item(27, 12).Value = "2025-05-25"
$ws.Cells.Item(27, 13).Value = "Scheduled"

# Row 28: MA261 - Differential Equations
$ws.Cells.Item(28, 1).Value = "MA261"
$ws.Cells.Item(28, 2).Value = "Differential Equations"
$ws.Cells.Item(28, 3).Value = "Theory"
$ws.Cells.Item(28, 4).Value = "2.0 hours"
$ws.Cells.Item(28, 5).Value = 120
$ws.Cells.Item(28, 6).Value = "CSE"
$ws.Cells.Item(28, 7).Value = 3
$ws.Cells.Item(28, 8).Value = "2025-11-28"
$ws.Cells.Item(28, 9).Value = "Friday"
$ws.Cells.Item(28, 10).Value = "Afternoon"
$ws.Cells.Item(28, 11).Value = "14:00 - 17:00"
$ws.Cells.Item(28, 12).Value = "2025-04-15"
$ws.Cells.Item(28, 13).Value = "Scheduled"

# Row 29: EC264 - Semiconductor Devices
$ws.Cells.Item(29, 1).Value = "EC264"
$ws.Cells.Item(29, 2).Value = "Semiconductor Devices"
$ws.Cells.Item(29, 3).Value = "Theory"
$ws.Cells.Item(29, 4).Value = "2.0 hours"
$ws.Cells.Item(29, 5).Value = 120
$ws.Cells.Item(29, 6).Value = "ECE"
$ws.Cells.Item(29, 7).Value = 5
$ws.Cells.Item(29, 8).Value = "2025-11-28"
$ws.Cells.Item(29, 9).Value = "Friday"
$ws.Cells.Item(29, 10).Value = "Afternoon"
$ws.Cells.Item(29, 11).Value = "14:00 - 17:00"
$ws.Cells.Item(29, 12).Value = "2025-05-26"
$ws.Cells.Item(29, 13).Value = "Scheduled"

# Row 30: MA262 - Multivariable Calculus
$ws.Cells.Item(30, 1).Value = "MA262"
$ws.Cells.Item(30, 2).Value = "Multivariable Calculus"
$ws.Cells.Item(30, 3).Value = "Theory"
$ws.Cells.Item(30, 4).Value = "2.0 hours"
$ws.Cells.Item(30, 5).Value = 120
$ws.Cells.Item(30, 6).Value = "CSE"
$ws.Cells.Item(30, 7).Value = 3
$ws.Cells.Item(30, 8).Value = "2025-12-01"
$ws.Cells.Item(30, 9).Value = "Monday"
$ws.Cells.Item(30, 10).Value = "Morning"
$ws.Cells.Item(30, 11).Value = "09:00 - 12:00"
$ws.Cells.Item(30, 12).Value = "2025-04-16"
$ws.Cells.Item(30, 13).Value = "Scheduled"

# Row 31: EC262 - Analog Electronics
$ws.Cells.Item(31, 1).Value = "EC262"
$ws.Cells.Item(31, 2).Value = "Analog Electronics"
$ws.Cells.Item(31, 3).Value = "Theory"
$ws.Cells.Item(31, 4).Value = "2.0 hours"
$ws.Cells.Item(31, 5).Value = 120
$ws.Cells.Item(31, 6).Value = "ECE"
$ws.Cells.Item(31, 7).Value = 5
$ws.Cells.Item(31, 8).Value = "2025-12-01"
$ws.Cells.Item(31, 9).Value = "Monday"
$ws.Cells.Item(31, 10).Value = "Morning"
$ws.Cells.Item(31, 11).Value = "09:00 - 12:00"
$ws.Cells.Item(31, 12).Value = "2025-05-27"
$ws.Cells.Item(31, 13).Value = "Scheduled"

# Row 32: CS261 - Operating Systems
$ws.Cells.Item(32, 1).Value = "CS261"
$ws.Cells.Item(32, 2).Value = "Operating Systems"
$ws.Cells.Item(32, 3).Value = "Theory"
$ws.Cells.Item(32, 4).Value = "3.0 hours"
$ws.Cells.Item(32, 5).Value = 180
$ws.Cells.Item(32, 6).Value = "CSE"
$ws.Cells.Item(32, 7).Value = 3
$ws.Cells.Item(32, 8).Value = "2025-12-01"
$ws.Cells.Item(32, 9).Value = "Monday"
$ws.Cells.Item(32, 10).Value = "Afternoon"
$ws.Cells.Item(32, 11).Value = "14:00 - 17:00"
$ws.Cells.Item(32, 12).Value = "2025-04-17"
$ws.Cells.Item(32, 13).Value = "Scheduled"

# Row 33: CS262 - Software Design
$ws.Cells.Item(33, 1).Value = "CS262"
$ws.Cells.Item(33, 2).Value = "Software Design"
$ws.Cells.Item(33, 3).Value = "Theory"
$ws.Cells.Item(33, 4).Value = "3.0 hours"
$ws.Cells.Item(33, 5).Value = 180
$ws.Cells.Item(33, 6).Value = "CSE"
$ws.Cells.Item(33, 7).Value = 3
$ws.Cells.Item(33, 8).Value = "2025-12-02"
$ws.Cells.Item(33, 9).Value = "Tuesday"
$ws.Cells.Item(33, 10).Value = "Morning"
$ws.Cells.Item(33, 11).Value = "09:00 - 12:00"
$ws.Cells.Item(33, 12).Value = "2025-04-18"
$ws.Cells.Item(33, 13).Value = "Scheduled"

# Row 34: CS253 - Introduction to AI
$ws.Cells.Item(34, 1).Value = "CS253"
$ws.Cells.Item(34, 2).Value = "Introduction to AI"
$ws.Cells.Item(34, 3).Value = "Theory"
$ws.Cells.Item(34, 4).Value = "2.0 hours"
$ws.Cells.Item(34, 5).Value = 120
$ws.Cells.Item(34, 6).Value = "CSE"
$ws.Cells.Item(34, 7).Value = 3
$ws.Cells.Item(34, 8).Value = "2025-12-02"
$ws.Cells.Item(34, 9).Value = "Tuesday"
$ws.Cells.Item(34, 10).Value = "Afternoon"
$ws.Cells.Item(34, 11).Value = "14:00 - 17:00"
$ws.Cells.Item(34, 12).Value = "2025-04-29"
$ws.Cells.Item(34, 13).Value = "Scheduled"

# Row 35: CS152 - Data Science with Python
$ws.Cells.Item(35, 1).Value = "CS152"
$ws.Cells.Item(35, 2).Value = "Data Science with Python"
$ws.Cells.Item(35, 3).Value = "Lab"
$ws.Cells.Item(35, 4).Value = "2.0 hours"
$ws.Cells.Item(35, 5).Value = 120
$ws.Cells.Item(35, 6).Value = "CSE"
$ws.Cells.Item(35, 7).Value = 3
$ws.Cells.Item(35, 8).Value = "2025-12-03"
$ws.Cells.Item(35, 9).Value = "Wednesday"
$ws.Cells.Item(35, 10).Value = "Morning"
$ws.Cells.Item(35, 11).Value = "09:00 - 12:00"
$ws.Cells.Item(35, 12).Value = "2025-04-30"
$ws.Cells.Item(35, 13).Value = "Scheduled"

# Row 36: CS309 - Statistics for CS
$ws.Cells.Item(36, 1).Value = "CS309"
$ws.Cells.Item(36, 2).Value = "Statistics for CS"
$ws.Cells.Item(36, 3).Value = "Theory"
$ws.Cells.Item(36, 4).Value = "3.0 hours"
$ws.Cells.Item(36, 5).Value = 180
$ws.Cells.Item(36, 6).Value = "CSE"
$ws.Cells.Item(36, 7).Value = 5
$ws.Cells.Item(36, 8).Value = "2025-12-03"
$ws.Cells.Item(36, 9).Value = "Wednesday"
$ws.Cells.Item(36, 10).Value = "Afternoon"
$ws.Cells.Item(36, 11).Value = "14:00 - 17:00"
$ws.Cells.Item(36, 12).Value = "2025-05-17"
$ws.Cells.Item(36, 13).Value = "Scheduled"

# Row 37: CS303 - Computer Networks
$ws.Cells.Item(37, 1).Value = "CS303"
$ws.Cells.Item(37, 2).Value = "Computer Networks"
$ws.Cells.Item(37, 3).Value = "Theory"
$ws.Cells.Item(37, 4).Value = "3.0 hours"
$ws.Cells.Item(37, 5).Value = 180
$ws.Cells.Item(37, 6).Value = "CSE"
$ws.Cells.Item(37, 7).Value = 5
$ws.Cells.Item(37, 8).Value = "2025-12-04"
$ws.Cells.Item(37, 9).Value = "Thursday"
$ws.Cells.Item(37, 10).Value = "Morning"
$ws.Cells.Item(37, 11).Value = "09:00 - 12:00"
$ws.Cells.Item(37, 12).Value = "2025-05-18"
$ws.Cells.Item(37, 13).Value = "Scheduled"

# Row 38: HS101 - Environmental Studies
$ws.Cells.Item(38, 1).Value = "HS101"
$ws.Cells.Item(38, 2).Value = "Environmental Studies"
$ws.Cells.Item(38, 3).Value = "Theory"
$ws.Cells.Item(38, 4).Value = "2.0 hours"
$ws.Cells.Item(38, 5).Value = 120
$ws.Cells.Item(38, 6).Value = "CSE"
$ws.Cells.Item(38, 7).Value = 5
$ws.Cells.Item(38, 8).Value = "2025-12-04"
$ws.Cells.Item(38, 9).Value = "Thursday"
$ws.Cells.Item(38, 10).Value = "Afternoon"
$ws.Cells.Item(38, 11).Value = "14:00 - 17:00"
$ws.Cells.Item(38, 12).Value = "2025-05-20"
$ws.Cells.Item(38, 13).Value = "Scheduled"

# Row 39: CS463 - Parallel Computing
$ws.Cells.Item(39, 1).Value = "CS463"
$ws.Cells.Item(39, 2).Value = "Parallel Computing"
$ws.Cells.Item(39, 3).Value = "Theory"
$ws.Cells.Item(39, 4).Value = "3.0 hours"
$ws.Cells.Item(39, 5).Value = 180
$ws.Cells.Item(39, 6).Value = "CSE"
$ws.Cells.Item(39, 7).Value = 5
$ws.Cells.Item(39, 8).Value = "2025-12-05"
$ws.Cells.Item(39, 9).Value = "Friday"
$ws.Cells.Item(39, 10).Value = "Morning"
$ws.Cells.Item(39, 11).Value = "09:00 - 12:00"
$ws.Cells.Item(39, 12).Value = "2025-05-29"
$ws.Cells.Item(39, 13).Value = "Scheduled"

# Row 40: CS308 - Compiler Design
$ws.Cells.Item(40, 1).Value = "CS308"
$ws.Cells.Item(40, 2).Value = "Compiler Design"
$ws.Cells.Item(40, 3).Value = "Theory"
$ws.Cells.Item(40, 4).Value = "3.0 hours"
$ws.Cells.Item(40, 5).Value = 180
$ws.Cells.Item(40, 6).Value = "CSE"
$ws.Cells.Item(40, 7).Value = 5
$ws.Cells.Item(40, 8).Value = "2025-12-05"
$ws.Cells.Item(40, 9).Value = "Friday"
$ws.Cells.Item(40, 10).Value = "Afternoon"
$ws.Cells.Item(40, 11).Value = "14:00 - 17:00"
$ws.Cells.Item(40, 12).Value = "2025-05-30"
$ws.Cells.Item(40, 13).Value = "Scheduled"

# Drop the temporary text format now that the values are safely stored as text.
$dateTextRangeH.ClearFormats()
$dateTextRangeL.ClearFormats()

